# Insert a new weekly price-report row at row 87 ("Hortaliza, Feria Lagunitas
# de Puerto Montt - Poroto verde"). All existing rows from 87 down get pushed
# down by one (Excel's normal row-insert semantics), and the freed row 87 is
# populated with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 87:132 down to 88:133, carrying formatting along (mirrors
# right-click > Insert in Excel).
$ws.Rows.Item(87).EntireRow.Insert()

# Populate the newly freed row 87 with the new observation.
$ws.Range("A87").Value = 4
$ws.Range("B87").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C87").Value = "Los Lagos"
$ws.Range("D87").Value = 45001
$ws.Range("E87").Value = 10
$ws.Range("F87").Value = 100112031
$ws.Range("G87").Value = "Poroto verde"
$ws.Range("H87").Value = "Magnum"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 30
$ws.Range("K87").Value = 35000
$ws.Range("L87").Value = 35000
$ws.Range("M87").Value = 35000
$ws.Range("N87").Value = '$/saco 25 kilos'
$ws.Range("O87").Value = "Región Metropolitana"
$ws.Range("P87").Value = 1400
$ws.Range("Q87").Value = 25
$ws.Range("R87").Value = "Hortaliza"
